# DIS-project - E_R.pptx edit
#
# The commit removes two now-unused, empty placeholder/text shapes from
# slide 2 ("User - database") and one empty placeholder shape from
# slide 3 ("Result - database"). Google Slides' exporter then re-numbers
# (compacts) every shape id/name that came after the removed shapes so
# there are no gaps left behind. We reproduce the same end result here:
# delete the now-empty shapes, then rename the shapes that shift "up"
# so the surviving shapes/names line up with the post-edit numbering.

$p = $ppt.ActivePresentation

# ----------------------------------------------------------------------
# Slide 2 ("User - database"): remove the empty body placeholder
# (old id 61) and the empty free textbox (old id 62), then renumber the
# remaining shape names down by two (63->61, 64->62, 65->63, 66->64,
# 67->65).
# ----------------------------------------------------------------------
$s2 = $p.Slides.Item(2)

# Delete the empty "body" placeholder (Google Shape;61;p14).
$s2.Shapes.Item(2).Delete()
# Delete the empty free textbox (Google Shape;62;p14).
$s2.Shapes.Item(2).Delete()

# What's left, in order: 63 (users), 64 (username), 65 (password),
# 66 (cxn), 67 (cxn) -> rename down to 61, 62, 63, 64, 65.
$s2.Shapes.Item(2).Name = "Google Shape;61;p14"
$s2.Shapes.Item(3).Name = "Google Shape;62;p14"
$s2.Shapes.Item(4).Name = "Google Shape;63;p14"
$s2.Shapes.Item(5).Name = "Google Shape;64;p14"
$s2.Shapes.Item(6).Name = "Google Shape;65;p14"

# ----------------------------------------------------------------------
# Slide 2's notes page: its shape names embed the (now shifted) slide
# numbering scheme too, so they move down by two (69->67, 70->68).
# ----------------------------------------------------------------------
$n2 = $s2.NotesPage
$n2.Shapes.Item(1).Name = "Google Shape;67;g2e4428c0439_0_5:notes"
$n2.Shapes.Item(2).Name = "Google Shape;68;g2e4428c0439_0_5:notes"

# ----------------------------------------------------------------------
# Slide 3 ("Result - database"): remove the empty body placeholder
# (old id 73), then renumber the remaining shape names down by three
# (74->71, 75->72, ..., 95->92).
# ----------------------------------------------------------------------
$s3 = $p.Slides.Item(3)

# Delete the empty "body" placeholder (Google Shape;73;p15).
$s3.Shapes.Item(2).Delete()

$names3 = @(
  "Google Shape;71;p15",
  "Google Shape;72;p15",
  "Google Shape;73;p15",
  "Google Shape;74;p15",
  "Google Shape;75;p15",
  "Google Shape;76;p15",
  "Google Shape;77;p15",
  "Google Shape;78;p15",
  "Google Shape;79;p15",
  "Google Shape;80;p15",
  "Google Shape;81;p15",
  "Google Shape;82;p15",
  "Google Shape;83;p15",
  "Google Shape;84;p15",
  "Google Shape;85;p15",
  "Google Shape;86;p15",
  "Google Shape;87;p15",
  "Google Shape;88;p15",
  "Google Shape;89;p15",
  "Google Shape;90;p15",
  "Google Shape;91;p15",
  "Google Shape;92;p15"
)

for ($i = 0; $i -lt $names3.Count; $i++) {
  $s3.Shapes.Item($i + 2).Name = $names3[$i]
}
